# Trade #26 closed at 2026-02-17 23:57:49 - unknown UNKNOWN +0.000%
#
# Updates the Summary metrics, the MarketMaking row on the Strategy Status
# sheet, and appends the newly closed trade (#26) to both the "All Trades"
# and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.93   # Current Capital
$summary.Range("B4").Value = 0.93      # Total P&L $
$summary.Range("B5").Value = 0.72      # Total P&L %
$summary.Range("B6").Value = 26        # Total Trades
$summary.Range("B7").Value = 15        # Winning Trades
$summary.Range("B9").Value = 57.69     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row, row 6) ---------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.93     # Capital
$status.Range("D6").Value = 26         # Trades
$status.Range("E6").Value = 0.93       # P&L $
$status.Range("F6").Value = 0.93       # P&L %
$status.Range("G6").Value = 57.69      # Win Rate %

# ---- New trade row (#26) appended to the trade logs -----------------------
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry, G Exit,
#          H Status, I P&L%, J P&L$, K Capital After, L Entry slip,
#          M Exit slip, N Confidence, O Entry reason, P Exit reason,
#          Q Duration (min)
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A27").Value = 26

    # Date/Time are stored as plain text in this workbook (not real Excel
    # dates), so force the cells to Text format before assigning values
    # that would otherwise be auto-recognized as a date/time literal.
    $ws.Range("B27").NumberFormat = "@"
    $ws.Range("B27").Value = "2026-02-17"
    $ws.Range("C27").Value = "23:57:42"

    $ws.Range("D27").Value = "MarketMaking"
    $ws.Range("E27").Value = "DOWN"
    $ws.Range("F27").Value = 0.22
    $ws.Range("G27").Value = 0.33
    $ws.Range("H27").Value = "CLOSED"
    $ws.Range("I27").Value = 50
    $ws.Range("J27").Value = 0.11
    $ws.Range("K27").Value = 100.93
    $ws.Range("L27").Value = 0
    $ws.Range("M27").Value = 0
    $ws.Range("N27").Value = 0.6
    $ws.Range("O27").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P27").Value = "early_exit"
    $ws.Range("Q27").Value = 0.14
}
